$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells are written as TEXT (matching the original inlineStr
# cell type) rather than being auto-coerced to numbers by COM for numeric-looking
# strings like "1.002". Force text format first, write values, then restore the
# default ('Normal') cell style so no stray number-format style lingers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$updates = @{
    'D2' = '30.480.99'
    'E2' = '  -0.91%  '
    'D3' = '1.893.84'
    'E3' = '  -0.72%  '
    'D4' = '1.002'
    'E4' = '  +0.10%  '
    'D5' = '238.96'
    'E5' = '  -0.43%  '
    'D6' = '1.002'
    'E6' = '  +0.08%  '
    'D7' = '0.4834'
    'E7' = '  -2.18%  '
    'D8' = '0.2900'
    'E8' = '  -2.13%  '
    'D9' = '0.06614'
    'E9' = '  -1.91%  '
    'D10' = '1.894.53'
    'E10' = '  -0.66%  '
    'D11' = '16.97'
    'E11' = '  -0.52%  '
    'D12' = '0.07395'
    'E12' = '  +0.53%  '
    'D13' = '5.194'
    'E13' = '  +0.54%  '
    'D14' = '89.25'
    'E14' = '  +0.97%  '
    'D15' = '0.6634'
    'E15' = '  -1.16%  '
    'D16' = '30.473.70'
    'E16' = '  -0.75%  '
    'D17' = '13.59'
    'E17' = '  +0.58%  '
    'D18' = '0.000007776'
    'E18' = '  -1.80%  '
    'D19' = '0.9996'
    'E19' = '  -0.11%  '
    'D20' = '2.154.66'
    'E20' = '  +0.00%  '
    'D21' = '5.408'
    'E21' = '  +1.76%  '
    'B22' = 'BinanceUSD'
    'C22' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D22' = '1.001'
    'E22' = '  +0.00%  '
    'B23' = 'BitcoinCash'
    'C23' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D23' = '221.56'
    'E23' = '  +13.68%  '
    'D24' = '6.206'
    'E24' = '  -0.98%  '
    'D25' = '9.410'
    'E25' = '  -2.20%  '
    'D26' = '163.84'
    'E26' = '  +0.40%  '
    'D27' = '18.66'
    'E27' = '  +0.36%  '
    'D28' = '1.945'
    'E28' = '  -0.20%  '
    'D29' = '1.441'
    'E29' = '  -2.19%  '
    'D30' = '4.344'
    'E30' = '  -2.49%  '
    'D31' = '0.09183'
    'E31' = '  +0.36%  '
    'D32' = '4.060'
    'E32' = '  +0.42%  '
    'D33' = '0.05082'
    'E33' = '  -3.36%  '
    'D34' = '0.7569'
    'E34' = '  +1.97%  '
    'D35' = '1.160'
    'E35' = '  +4.49%  '
    'D36' = '2.706'
    'E36' = '  -0.82%  '
    'D37' = '0.01889'
    'E37' = '  +3.52%  '
    'D38' = '2.659'
    'E38' = '  -2.02%  '
    'D39' = '2.103'
    'E39' = '  +0.99%  '
    'D40' = '0.9201'
    'E40' = '  -0.20%  '
    'D41' = '6.031'
    'E41' = '  +1.30%  '
    'D42' = '107.23'
    'E42' = '  +0.40%  '
    'D43' = '0.4352'
    'E43' = '  -2.08%  '
    'D44' = '1.003'
    'D45' = '7.636'
    'E45' = '  +1.15%  '
    'B46' = 'NEARProtocol'
    'C46' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D46' = '1.609'
    'E46' = '  +12.37%  '
    'B47' = 'Algorand'
    'C47' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D47' = '0.1343'
    'E47' = '  -3.27%  '
    'D48' = '65.24'
    'E48' = '  -12.49%  '
    'D49' = '8.939'
    'E49' = '  -1.42%  '
    'D50' = '34.44'
    'E50' = '  -2.90%  '
    'D51' = '0.05707'
    'E51' = '  -2.58%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# Restore default styling (remove the temporary text-number-format override).
$dataRange.Style = "Normal"
